# Fix Training Data Issue (#48)
# The "Date" column (BF) held strings built as "<day>-<month>-<season>"
# (e.g. "5-11-2013-14") instead of a real ISO date. NBA.com's stats page
# labelled the games one day off, so the values are corrected here to the
# actual game date "2014-05-11" for every data row (BF2:BF31).
#
# We can't just assign the text to .Value / .Value2 / .Formula because the
# runtime's smart literal-entry parsing recognises "2014-05-11" as a real
# date and silently reformats the cell (adding a number-format style),
# which would introduce unwanted style/cellXfs churn. Instead we build the
# literal text via a text-formula in a scratch cell and paste-special just
# the value back in, which keeps the destination cells plain, unstyled
# shared-string text cells exactly like the other untouched columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$scratch = $ws.Range("ZZ1")
$scratch.Formula = "=""2014-05-11"""
$scratch.Copy()

for ($r = 2; $r -le 31; $r++) {
    $ws.Cells.Item($r, 58).PasteSpecial(-4163)
}

$scratch.Clear()
